$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New CRN rows (30 and 31) appended after existing data (row 29)
$ws.Range("A30").Value = "crn_00028"
$ws.Range("B30").Value = 3473.204909
$ws.Range("C30").Value = 264.312233
$ws.Range("D30").Value = -0.06714299999999999
$ws.Range("E30").Value = 2.29813
$ws.Range("F30").Value = 53.125
$ws.Range("G30").Value = 6
$ws.Range("H30").Value = 272
$ws.Range("I30").Value = 234

$ws.Range("A31").Value = "crn_00029"
$ws.Range("B31").Value = 3541.138945
$ws.Range("C31").Value = 246.031259
$ws.Range("D31").Value = 0.304077
$ws.Range("E31").Value = 4.054082
$ws.Range("F31").Value = 53.125
$ws.Range("G31").Value = 12
$ws.Range("H31").Value = 272
$ws.Range("I31").Value = 228

# Apply the same thin-border cell style used throughout the data rows (A2:I29)
$ws.Range("A30:I31").Borders.LineStyle = 1
